# Update cryptocurrency price (column D) and 1h volume change (column E)
# values on the active worksheet to reflect the latest scrape, as
# captured in the GitHub Actions commit "Updated cryptos list".
#
# Price-column values are stored as literal text in the source data
# (e.g. "3.198.73" uses "." as a thousands separator, and values like
# "217.74" must stay text, not be reinterpreted as numbers). A leading
# apostrophe forces Excel to keep an otherwise numeric-looking entry as
# text, matching the original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.949.60"
$ws.Range("E2").Value = "  +2.93%  "
$ws.Range("D3").Value = "3.198.73"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'217.74"
$ws.Range("E5").Value = "  +7.17%  "
$ws.Range("D6").Value = "'647.79"
$ws.Range("E6").Value = "  +7.06%  "
$ws.Range("D7").Value = "'0.394"
$ws.Range("E7").Value = "  +4.87%  "
$ws.Range("D8").Value = "'0.692"
$ws.Range("E8").Value = "  +5.42%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.196.58"
$ws.Range("D11").Value = "'0.574"
$ws.Range("E11").Value = "  +8.59%  "
$ws.Range("D12").Value = "'0.180"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  +5.98%  "
$ws.Range("D14").Value = "'5.40"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").Value = "'33.40"
$ws.Range("E15").Value = "  +5.11%  "
$ws.Range("D16").Value = "3.781.60"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "89.636.45"
$ws.Range("E17").Value = "  +3.03%  "
$ws.Range("D18").Value = "3.193.65"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'3.43"
$ws.Range("E19").Value = "  +15.06%  "
$ws.Range("D20").Value = "'0.0000225"
$ws.Range("E20").Value = "  +75.74%  "
$ws.Range("D21").Value = "'13.54"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").Value = "'437.91"
$ws.Range("E22").Value = "  +6.35%  "
$ws.Range("D23").Value = "'8.66"
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("D24").Value = "'5.10"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "'5.31"
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("D26").Value = "'11.96"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").Value = "'81.63"
$ws.Range("E27").Value = "  +11.48%  "
$ws.Range("D28").Value = "3.360.53"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'548.72"
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("D33").Value = "'4.05"
$ws.Range("E33").Value = "  +36.54%  "
$ws.Range("D34").Value = "'8.49"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("D35").Value = "'7.07"
$ws.Range("E35").Value = "  +6.88%  "
$ws.Range("D36").Value = "'1.94"
$ws.Range("E36").Value = "  +6.41%  "
$ws.Range("D37").Value = "'1.32"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").Value = "'22.54"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("D44").Value = "'0.376"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").Value = "'146.03"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "'174.38"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").Value = "'44.04"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("E48").Value = "  +9.99%  "
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'0.623"
$ws.Range("E51").Value = "  +7.06%  "
